$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "git" task row (row 18). Rows below it (the "ubuntu shutdown"
# row) shift up by one, so the former row 19 becomes the new row 18.
$ws.Rows("18:18").Delete()

# Re-touch the sheet's final row so it stays materialized after the shift
# (no-op toggle: leaves it exactly as it was, just keeps it present).
$ws.Rows(1048576).Hidden = $true
$ws.Rows(1048576).Hidden = $false

# Select the cell where the cursor lands after deleting the row.
$ws.Range("A18").Select()
